# Update market-price / profit figures across the Leve profit sheets
# (refreshed figures from the scheduled market-data runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1405.4375
$ws.Range("I70").Value = 1318.7
$ws.Range("K70").Value = 3956.1
$ws.Range("M70").Value = -3686.1
$ws.Range("H73").Value = 1405.4375
$ws.Range("I73").Value = 1318.7
$ws.Range("K73").Value = 3956.1
$ws.Range("M73").Value = -3020.1
$ws.Range("H74").Value = 4233.2666
$ws.Range("I74").Value = 2900
$ws.Range("K74").Value = 2900
$ws.Range("M74").Value = -1964
$ws.Range("H77").Value = 4233.2666
$ws.Range("I77").Value = 2900
$ws.Range("K77").Value = 14500
$ws.Range("M77").Value = -9820
$ws.Range("H113").Value = 43481820
$ws.Range("I113").Value = 55558188
$ws.Range("K113").Value = 55558188
$ws.Range("M113").Value = -55554934
$ws.Range("H116").Value = 4921.6665
$ws.Range("I116").Value = 4500
$ws.Range("J116").Value = 4974.375
$ws.Range("K116").Value = 4500
$ws.Range("L116").Value = 4974.375
$ws.Range("M116").Value = -1058
$ws.Range("N116").Value = -11858.375
$ws.Range("H129").Value = 228224.1
$ws.Range("J129").Value = 286853.78
$ws.Range("L129").Value = 860561.3400000001
$ws.Range("N129").Value = -870561.3400000001
$ws.Range("H132").Value = 2105.4717
$ws.Range("I132").Value = 2266.1777
$ws.Range("J132").Value = 1201.5
$ws.Range("K132").Value = 6798.533100000001
$ws.Range("L132").Value = 3604.5
$ws.Range("M132").Value = -4268.533100000001
$ws.Range("N132").Value = -8664.5
$ws.Range("H135").Value = 29415712
$ws.Range("I135").Value = 1256.25
$ws.Range("K135").Value = 11306.25
$ws.Range("M135").Value = -8771.25
$ws.Range("H138").Value = 28573928
$ws.Range("I138").Value = 58825156
$ws.Range("J138").Value = 3320.7778
$ws.Range("K138").Value = 176475468
$ws.Range("L138").Value = 9962.3334
$ws.Range("M138").Value = -176470328
$ws.Range("N138").Value = -20242.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1404.45
$ws.Range("I61").Value = 1344.1578
$ws.Range("J61").Value = 2550
$ws.Range("K61").Value = 1344.1578
$ws.Range("L61").Value = 2550
$ws.Range("M61").Value = -1132.1578
$ws.Range("N61").Value = -2974
$ws.Range("H74").Value = 50001884
$ws.Range("I74").Value = 90909700
$ws.Range("J74").Value = 3443.6667
$ws.Range("K74").Value = 90909700
$ws.Range("L74").Value = 3443.6667
$ws.Range("M74").Value = -90908826
$ws.Range("N74").Value = -5191.6667
$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20676
$ws.Range("H77").Value = 50001884
$ws.Range("I77").Value = 90909700
$ws.Range("J77").Value = 3443.6667
$ws.Range("K77").Value = 454548500
$ws.Range("L77").Value = 17218.3335
$ws.Range("M77").Value = -454544132
$ws.Range("N77").Value = -25954.3335
$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22340
$ws.Range("H97").Value = 71429600
$ws.Range("I97").Value = 772.7273
$ws.Range("J97").Value = 333335300
$ws.Range("K97").Value = 772.7273
$ws.Range("L97").Value = 333335300
$ws.Range("M97").Value = -276.7273
$ws.Range("N97").Value = -333336292
$ws.Range("H110").Value = 640.5
$ws.Range("I110").Value = 537.5833
$ws.Range("J110").Value = 949.25
$ws.Range("K110").Value = 537.5833
$ws.Range("L110").Value = 949.25
$ws.Range("M110").Value = 1507.4167
$ws.Range("N110").Value = -5039.25
$ws.Range("H136").Value = 1404.45
$ws.Range("I136").Value = 1344.1578
$ws.Range("J136").Value = 2550
$ws.Range("K136").Value = 4032.4734
$ws.Range("L136").Value = 7650
$ws.Range("M136").Value = -1482.4734
$ws.Range("N136").Value = -12750

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1151.125
$ws.Range("I94").Value = 952.25
$ws.Range("K94").Value = 952.25
$ws.Range("M94").Value = -501.25
$ws.Range("H107").Value = 852.625
$ws.Range("I107").Value = 965.8333
$ws.Range("K107").Value = 965.8333
$ws.Range("M107").Value = 954.1667
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4232.7666
$ws.Range("I31").Value = 5487.3335
$ws.Range("J31").Value = 3919.125
$ws.Range("K31").Value = 5487.3335
$ws.Range("L31").Value = 3919.125
$ws.Range("M31").Value = -5192.3335
$ws.Range("N31").Value = -4509.125
$ws.Range("H34").Value = 4232.7666
$ws.Range("I34").Value = 5487.3335
$ws.Range("J34").Value = 3919.125
$ws.Range("K34").Value = 5487.3335
$ws.Range("L34").Value = 3919.125
$ws.Range("M34").Value = -5285.3335
$ws.Range("N34").Value = -4323.125
$ws.Range("H62").Value = 71432480
$ws.Range("I62").Value = 142861340
$ws.Range("J62").Value = 3614.2856
$ws.Range("K62").Value = 142861340
$ws.Range("L62").Value = 3614.2856
$ws.Range("M62").Value = -142860716
$ws.Range("N62").Value = -4862.2856
$ws.Range("H65").Value = 71432480
$ws.Range("I65").Value = 142861340
$ws.Range("J65").Value = 3614.2856
$ws.Range("K65").Value = 714306700
$ws.Range("L65").Value = 18071.428
$ws.Range("M65").Value = -714303580
$ws.Range("N65").Value = -24311.428

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 573.25
$ws.Range("J98").Value = 395
$ws.Range("L98").Value = 1185
$ws.Range("N98").Value = -4181
$ws.Range("H131").Value = 712.88
$ws.Range("J131").Value = 712.88
$ws.Range("L131").Value = 2138.64
$ws.Range("N131").Value = -12218.64

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3946.6843
$ws.Range("I80").Value = 3120.8333
$ws.Range("J80").Value = 4327.846
$ws.Range("K80").Value = 3120.8333
$ws.Range("L80").Value = 4327.846
$ws.Range("M80").Value = -2122.8333
$ws.Range("N80").Value = -6323.846
$ws.Range("H83").Value = 3946.6843
$ws.Range("I83").Value = 3120.8333
$ws.Range("J83").Value = 4327.846
$ws.Range("K83").Value = 15604.1665
$ws.Range("L83").Value = 21639.23
$ws.Range("M83").Value = -10612.1665
$ws.Range("N83").Value = -31623.23
$ws.Range("H97").Value = 1030.7037
$ws.Range("I97").Value = 944.5
$ws.Range("J97").Value = 1410
$ws.Range("K97").Value = 944.5
$ws.Range("L97").Value = 1410
$ws.Range("M97").Value = -448.5
$ws.Range("N97").Value = -2402
$ws.Range("H102").Value = 15153763
$ws.Range("I102").Value = 17243444
$ws.Range("K102").Value = 17243444
$ws.Range("M102").Value = -17241822
$ws.Range("H126").Value = 4172.222
$ws.Range("I126").Value = 2690
$ws.Range("J126").Value = 6025
$ws.Range("K126").Value = 8070
$ws.Range("L126").Value = 18075
$ws.Range("M126").Value = -5600
$ws.Range("N126").Value = -23015

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2410.4583
$ws.Range("I7").Value = 2313.5
$ws.Range("J7").Value = 2701.3333
$ws.Range("K7").Value = 2313.5
$ws.Range("L7").Value = 2701.3333
$ws.Range("M7").Value = -2201.5
$ws.Range("N7").Value = -2925.3333
$ws.Range("H126").Value = 2410.4583
$ws.Range("I126").Value = 2313.5
$ws.Range("J126").Value = 2701.3333
$ws.Range("K126").Value = 6940.5
$ws.Range("L126").Value = 8103.999899999999
$ws.Range("M126").Value = -4470.5
$ws.Range("N126").Value = -13043.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1288.8889
$ws.Range("I96").Value = 1375
$ws.Range("J96").Value = 1220
$ws.Range("K96").Value = 1375
$ws.Range("L96").Value = 1220
$ws.Range("M96").Value = -2
$ws.Range("N96").Value = -3966
$ws.Range("H126").Value = 1332.8387
$ws.Range("I126").Value = 1211.875
$ws.Range("K126").Value = 3635.625
$ws.Range("M126").Value = -1165.625
$ws.Range("H139").Value = 50884.75
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 50884.75
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 50884.75
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -61164.75

